$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the _GoBack bookmark that currently sits after the very first
#    "Define the problem" paragraph (Problem 1 - A Cat, a Parrot...).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. "Socks in the Dark" section (Problem 2) - flesh out the
#    "Define the problem" sub-bullets.
#    Paragraph 20 is the empty ilvl=1 bullet right after "Define the problem:"
# ---------------------------------------------------------------------------
$pDefine = $d.Paragraphs.Item(20)
$pDefine.Range.InsertAfter("There are different amount of pairs of three different colors and I have to find out how many I have to pull out of my drawer to get matching socks (something I do everyday).")

$pDefine = $d.Paragraphs.Item(20)
$pDefine.Range.InsertParagraphAfter()
$pInsight = $d.Paragraphs.Item(21)
$pInsight.Range.InsertAfter("Some insight? There’s going to be a lot of sock pulling.")

$pInsight.Range.InsertParagraphAfter()
$pGoal = $d.Paragraphs.Item(22)
$pGoal.Range.InsertAfter("The overall goal is to get a pair of each of the colors of socks while not looking at the socks.")

# ---------------------------------------------------------------------------
# 3. "Break the problem apart" sub-bullets. After the inserts above the old
#    "Identify potential solutions" bullet (originally ilvl=0) has shifted
#    down and becomes the first ilvl=1 sub-bullet under "Break the problem
#    apart:" - demote it and rewrite its text to the constraints sentence.
# ---------------------------------------------------------------------------
$pConstraints = $d.Paragraphs.Item(24)
$pConstraints.Range.Text = "The constraints are that you can’t look at the socks and that you are pulling one sock at a time."
$pConstraints = $d.Paragraphs.Item(24)
$pConstraints.Range.ListFormat.ListLevelNumber = 2

$pConstraints.Range.InsertParagraphAfter()
$pSubGoals = $d.Paragraphs.Item(25)
$pSubGoals.Range.InsertAfter("The sub-goals are to be able to suddenly see the color of the sock even though you are picking them out in the dark, and to be able to pick out socks in pairs.")

# ---------------------------------------------------------------------------
# 4. The old "Evaluate potential solutions" bullet becomes the (real)
#    "Identify potential solutions" bullet, and it carries the relocated
#    _GoBack bookmark right at the end of its text, inside the paragraph.
#    A brand-new paragraph is then added below it with the actual
#    "Evaluate potential solutions" text.
# ---------------------------------------------------------------------------
$pIdentify = $d.Paragraphs.Item(26)
# Temporarily append a placeholder character so that the insertion point for
# the bookmark is not the literal last character of the paragraph (adding a
# bookmark collapsed exactly at paragraph-end-minus-one position is mishandled
# by this host) - add the bookmark next to the placeholder, then delete it.
$pIdentify.Range.Text = "Identify potential solutionsX"
$pIdentify = $d.Paragraphs.Item(26)
$bmPos = $pIdentify.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$placeholderRange = $d.Range($bmPos, $bmPos + 1)
$placeholderRange.Delete()

$pIdentify = $d.Paragraphs.Item(26)
$pIdentify.Range.InsertParagraphAfter()
$pEvaluate = $d.Paragraphs.Item(27)
$pEvaluate.Range.InsertAfter("Evaluate potential solutions")
